$wb = $excel.ActiveWorkbook

$oldGuid = "69cb68e0-873b-4c03-96db-f069882655c1"
$newGuid = "f377ad96-056b-4751-bba6-cdaca01d82cb"

$oldZhHash = "706a8c15fa8a40a9fd6880256ade643bdef32427"
$newZhHash = "b1922a8510236da5d28bf84a0494413f8576a390"
$oldDeHash = "706a8c15fa8a40a9fd6880256ade643bdef32427"
$newDeHash = "b1922a8510236da5d28bf84a0494413f8576a390"

$newHoGenDate = "2017-02-17 09:02:16"
$newZhHandoffDate = "2017-02-17 09:02:00"

# --- Sheet "Overview" ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "e2e\$newGuid.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$ws.Range("G2").Value = $newHoGenDate

# --- Sheet "zh-cn" ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$ws.Range("H2").Value = $newZhHandoffDate

# --- Sheet "de-de" ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"
$ws.Range("H2").Value = $newHoGenDate
